$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.982.83"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "'1.848.59"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'309.93"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").Value = "'0.3679"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").Value = "'0.07228"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "'0.9288"
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "'0.07744"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "'1.858.91"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "'5.342"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "'6.440"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "'88.79"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "'0.000008650"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "'27.003.28"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'14.47"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").Value = "'5.063"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("D24").Value = "'1.931"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").Value = "'152.85"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'18.23"
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("D27").Value = "'2.003"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'114.24"
$ws.Range("D29").Value = "'4.963"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").Value = "'0.08879"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "'3.324"
$ws.Range("E31").Value = "  +5.37%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'0.7428"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").Value = "'4.505"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").Value = "'2.746"
$ws.Range("E35").Value = "  -3.33%  "
$ws.Range("D36").Value = "'1.114"
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "'0.05267"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("D39").Value = "'2.980"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").Value = "'0.5218"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("D41").Value = "'6.986"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").Value = "'8.230"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").Value = "'10.60"
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("D45").Value = "'0.4736"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "'1.013"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'101.76"
$ws.Range("E47").Value = "  +3.49%  "
$ws.Range("D48").Value = "'1.609"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "'65.75"
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'0.8878"
$ws.Range("E51").Value = "  +4.16%  "
